$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 75 and 76 (appended after existing row 74)
$rowsData = @(
    @{ Row = 75; A = 45505; B = 830.9181992222; C = 222.11206359; D = 0; E = 0; F = 0; G = 0; I = 245.5479830475; J = 0; K = 0.054048930736; L = 0; M = 0; N = 127.23763173696; O = 58.43031527700001; P = 0; Q = 0.0000025224; R = 0; S = 0; T = 0; U = 366.3887079303344; W = 0; X = 0; Y = 0; Z = 255.9178744242 },
    @{ Row = 76; A = 45506; B = 781.8965324363; C = 207.2886453235; D = 0; E = 0; F = 0; G = 0; I = 224.1344537989; J = 0; K = 0.053085160524; L = 0; M = 0; N = 114.64117319744; O = 55.158786583; P = 0; Q = 0.0000022368; R = 0; S = 0; T = 0; U = 326.8586413275155; W = 0; X = 0; Y = 0; Z = 237.629109251934 }
)

$cols = @("B","C","D","E","F","G","I","J","K","L","M","N","O","P","Q","R","S","T","U","W","X","Y","Z")

foreach ($rd in $rowsData) {
    $r = $rd.Row

    # Copy column A's cell (date, style s="2") from the prior row so the
    # new date cell keeps the same formatting as the rest of the column.
    $ws.Range("A74").Copy($ws.Range("A$r"))
    $ws.Range("A$r").Value = $rd["A"]

    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = $rd[$c]
    }
}
